$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date (column C) for rows 2-6 from 2023-09-01 (45170) to 2023-09-05 (45174)
foreach ($r in 2..6) {
    $ws.Cells.Item($r, 3).Value = 45174
}
